$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP")
if (-not $ws) { $ws = $wb.ActiveSheet }

# New formulas for O5/O6: replace the old SQRT bound with a FLOOR.MATH based
# integer bound (the "fix" being made in this commit).
$formula = "=_xlfn.FLOOR.MATH((-1+(1-4*(-2)*`$M`$8)^0.5)/2)"
$ws.Range("O5").Formula = $formula
$ws.Range("O6").Formula = $formula

# New column P: labels "int" next to the two bound cells, formatted with the
# built-in Comma number format (numFmtId 43).
$ws.Range("P5").Value = "int"
$ws.Range("P6").Value = "int"
$ws.Range("P5:P6").NumberFormat = "_(* #,##0.00_);_(* (#,##0.00);_(* ""-""??_);_(@_)"

# Size the new column to fit its contents, like Excel would do automatically
# when a user types into a fresh column.
$ws.Columns("P:P").AutoFit() | Out-Null
